$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the report forward by one month ---------------------------------
# The "closing" figures of the previous month become this month's updated
# numbers (new month data appended). Update every static value in column B
# to the values for the new month.

$ws.Range("B1").Value = 45413
$ws.Range("B3").Value = 173900
$ws.Range("B4").Value = 141125
$ws.Range("B7").Value = 69550
$ws.Range("B10").Value = 69550
$ws.Range("B11").Value = 54225
$ws.Range("B12").Value = 191250
$ws.Range("B14").Value = 2456
$ws.Range("B15").Value = 1254.9000000000001
$ws.Range("B17").Value = 879.8
$ws.Range("B20").Value = 879.8
$ws.Range("B21").Value = 205
$ws.Range("B22").Value = 2626.1000000000004
$ws.Range("B24").Value = 45
$ws.Range("B27").ClearContents()
$ws.Range("B30").Value = 0
$ws.Range("B33").Value = 193921.1
$ws.Range("B35").Value = 70429.8
$ws.Range("B38").Value = 70429.8
$ws.Range("B40").Value = 63524
$ws.Range("B41").Value = 11150
$ws.Range("B43").Value = 74674
$ws.Range("B44").Value = 6905.8000000000029
$ws.Range("B45").Value = 0.098052244930413021
$ws.Range("B48").Value = 43701.899999999994
$ws.Range("B49").Value = 63524
$ws.Range("B50").Value = 11150
$ws.Range("B51").Value = 651.5
$ws.Range("B52").Value = 119027.4
$ws.Range("B53").Value = 62357.8
$ws.Range("B54").Value = 9050
$ws.Range("B56").Value = 71407.8
$ws.Range("B57").Value = 47619.599999999991
$ws.Range("B59").Value = 63142.9
$ws.Range("B60").Value = 62357.8
$ws.Range("B61").Value = 125500.70000000001
$ws.Range("B63").Value = 61041.3
$ws.Range("B64").Value = 651.5
$ws.Range("B65").Value = 98.2
$ws.Range("B66").Value = 6100.4
$ws.Range("B68").Value = 67891.399999999994
$ws.Range("B70").Value = 57609.300000000017
$ws.Range("B72").Value = 105228.90000000001
$ws.Range("B75").Value = 286264.5
$ws.Range("B76").Value = 61041.3
$ws.Range("B77").Value = 13264.5
$ws.Range("B79").Value = 96810.9
$ws.Range("B83").Value = 263759.40000000002
$ws.Range("B85").Value = 1576.0999999999995
$ws.Range("B86").Value = 298.3
$ws.Range("B88").Value = 69.5
$ws.Range("B89").Value = 1130.5
$ws.Range("B91").ClearContents()
$ws.Range("B92").Value = 5918
$ws.Range("B93").Value = 7416.3
$ws.Range("B95").Value = 83.6
$ws.Range("B96").Value = 298.3
$ws.Range("B97").Value = 949.4
$ws.Range("B99").ClearContents()
$ws.Range("B100").Value = 6605
$ws.Range("B101").Value = 7936.3
$ws.Range("B102").Value = 1056.0999999999995
$ws.Range("B104").Value = 370044.4
$ws.Range("B106").Value = 2156.5
$ws.Range("B107").Value = 13104.400000000001
$ws.Range("B110").Value = 4134
$ws.Range("B111").Value = 11126.900000000001
$ws.Range("B113").Value = 571086.39999999991
$ws.Range("B114").Value = 163060.69999999998
$ws.Range("B115").Value = 734147.09999999986
$ws.Range("B116").Value = 563965.5
$ws.Range("B117").Value = 170181.59999999986
$ws.Range("B118").Value = 168022.19999999998
$ws.Range("B119").Value = 2159.3999999998778
$ws.Range("B120").Value = 0.012688798318971496
$ws.Range("B122").Value = 176401
$ws.Range("B123").Value = 163060.69999999998
$ws.Range("B124").Value = 339461.69999999995
$ws.Range("B125").Value = 193921.1
$ws.Range("B126").Value = 145540.59999999995
$ws.Range("B127").Value = 24640.999999999942
$ws.Range("B128").Value = 170181.59999999989
$ws.Range("B129").Value = 168022.19999999998
$ws.Range("B130").Value = 2159.3999999999069
$ws.Range("B131").Value = 0.012688798318971664

# --- Unmerge the section-header rows (A2:B2, A13:B13, ...) -----------------
$ws.Range("A2:B2").UnMerge()
$ws.Range("A13:B13").UnMerge()
$ws.Range("A23:B23").UnMerge()
$ws.Range("A34:B34").UnMerge()
$ws.Range("A47:B47").UnMerge()
$ws.Range("A58:B58").UnMerge()
$ws.Range("A74:B74").UnMerge()
$ws.Range("A84:B84").UnMerge()
$ws.Range("A105:B105").UnMerge()
$ws.Range("A112:B112").UnMerge()
$ws.Range("A121:B121").UnMerge()

# A handful of header rows pick up an explicit (but unchanged) row height
# once they're unmerged.
$ws.Rows.Item(34).RowHeight = 15.75
$ws.Rows.Item(47).RowHeight = 15.75
$ws.Rows.Item(58).RowHeight = 15.75
$ws.Rows.Item(84).RowHeight = 15.75

# --- Update the view: reset scroll position, zoom out, and move the
#     selection to B1 -------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("B1").Select() | Out-Null
